# Gantt Chart updates (per commit: "edited the flow chart")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt Chart")

# Row 8 - Gantt Chart Write Up: progress 63% -> 65%
$ws.Range("I8").Value = 0.65

# Row 10 - Flow Chart: end date pushed out, days recalcs automatically
$ws.Range("F10").Value = 45737

# Row 13 - Creating 5 Levels(Maps): status -> Complete
$ws.Range("H13").Value = "Complete"

# Row 15 - 2 Enemy Character Dev: status -> Complete, 100% done
$ws.Range("H15").Value = "Complete"
$ws.Range("I15").Value = 1

# Row 16 - Champ. Player Special Abilities Dev: status -> In progress
$ws.Range("H16").Value = "In progress"

# Row 17 - renamed from "Bomb Development" to "Bomb animation"; dates shifted
$ws.Range("C17").Value = "Bomb animation"
$ws.Range("E17").Value = 45737
$ws.Range("F17").Value = 45744

# Row 18 - 2 Enemy's Special Abilities Dev: end date shifted, in progress, 80%
$ws.Range("F18").Value = 45737
$ws.Range("H18").Value = "In progress"
$ws.Range("I18").Value = 0.8

# Row 19 - 2 Enemy's Special Abilities Dev: end date shifted, in progress, 50%
$ws.Range("F19").Value = 45737
$ws.Range("H19").Value = "In progress"
$ws.Range("I19").Value = 0.5

# Row 20 - renamed from "The Special Door Dev" to "Bomb destruction"; dates shifted
$ws.Range("C20").Value = "Bomb destruction"
$ws.Range("E20").Value = 45744
$ws.Range("F20").Value = 45751

# Row 22 - renamed from "Level-Up Dev" to "Comment Background and Tile Manager Class"
$ws.Range("C22").Value = "Comment Background and Tile Manager Class"
$ws.Range("E22").Value = 45733
$ws.Range("F22").Value = 45737
$ws.Range("H22").Value = "In progress"

# Row 23 - renamed from "Champ. Player Health Dev" to "Special door Dev"; reassigned; dates shifted
$ws.Range("C23").Value = "Special door Dev"
$ws.Range("D23").Value = "Andrias"
$ws.Range("E23").Value = 45737
$ws.Range("F23").Value = 45744

# Row 24 - renamed from "Game Timer Dev" to "Comment all Enemy Classes"; reassigned; dates shifted
$ws.Range("C24").Value = "Comment all Enemy Classes"
$ws.Range("D24").Value = "Both"
$ws.Range("E24").Value = 45737
$ws.Range("F24").Value = 45744

# Row 25 - renamed from "Player Login GUI" to "Comment Character and JackBomber Classes"; reassigned; dates shifted
$ws.Range("C25").Value = "Comment Character and JackBomber Classes"
$ws.Range("D25").Value = "Murat C. GZ"
$ws.Range("E25").Value = 45737
$ws.Range("F25").Value = 45744

# Rows 26-29 - activity names removed (tasks dropped from the plan)
$ws.Range("C26").Value = ""
$ws.Range("C27").Value = ""
$ws.Range("C28").Value = ""
$ws.Range("C29").Value = ""
